# Scheduled runner update: refresh cached market-board profit calculations
# (currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns H, I, J, K, L, M, N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and
# WVR sheets with newly sampled prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2028.4166
$ws.Range("J9").Value = 1355.3334
$ws.Range("L9").Value = 1355.3334
$ws.Range("N9").Value = -1693.3334

$ws.Range("H19").Value = 1854.2916
$ws.Range("I19").Value = 1185.2858
$ws.Range("J19").Value = 2129.7646
$ws.Range("K19").Value = 1185.2858
$ws.Range("L19").Value = 2129.7646
$ws.Range("M19").Value = -1010.2858
$ws.Range("N19").Value = -2479.7646

$ws.Range("H96").Value = 66667680
$ws.Range("I96").Value = 976
$ws.Range("J96").Value = 200001090
$ws.Range("K96").Value = 2928
$ws.Range("L96").Value = 600003270
$ws.Range("M96").Value = -1555
$ws.Range("N96").Value = -600006016

$ws.Range("H100").Value = 21332.5
$ws.Range("I100").Value = 4999
$ws.Range("J100").Value = 24599.2
$ws.Range("K100").Value = 4999
$ws.Range("L100").Value = 24599.2
$ws.Range("M100").Value = -4458
$ws.Range("N100").Value = -25681.2

$ws.Range("H116").Value = 5563331.5
$ws.Range("I116").Value = 7943059.5
$ws.Range("J116").Value = 10632.833
$ws.Range("K116").Value = 7943059.5
$ws.Range("L116").Value = 10632.833
$ws.Range("M116").Value = -7939617.5
$ws.Range("N116").Value = -17516.833

$ws.Range("H137").Value = 6518.7095
$ws.Range("I137").Value = 3356.353
$ws.Range("J137").Value = 10358.714
$ws.Range("K137").Value = 10069.059
$ws.Range("L137").Value = 31076.142
$ws.Range("M137").Value = -7519.059000000001
$ws.Range("N137").Value = -36176.142

$ws.Range("H138").Value = 5850.2246
$ws.Range("I138").Value = 2633.5833
$ws.Range("J138").Value = 6893.4595
$ws.Range("K138").Value = 7900.749899999999
$ws.Range("L138").Value = 20680.3785
$ws.Range("M138").Value = -2760.749899999999
$ws.Range("N138").Value = -30960.3785

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8304.895
$ws.Range("I2").Value = 6020.9287
$ws.Range("J2").Value = 14700
$ws.Range("K2").Value = 6020.9287
$ws.Range("L2").Value = 14700
$ws.Range("M2").Value = -5907.9287
$ws.Range("N2").Value = -14926

$ws.Range("H61").Value = 5205.3716
$ws.Range("I61").Value = 4765.3105
$ws.Range("K61").Value = 4765.3105
$ws.Range("M61").Value = -4553.3105

$ws.Range("H74").Value = 40003344
$ws.Range("I74").Value = 71431150
$ws.Range("K74").Value = 71431150
$ws.Range("M74").Value = -71430276

$ws.Range("H77").Value = 40003344
$ws.Range("I77").Value = 71431150
$ws.Range("K77").Value = 357155750
$ws.Range("M77").Value = -357151382

$ws.Range("H116").Value = 8304.895
$ws.Range("I116").Value = 6020.9287
$ws.Range("J116").Value = 14700
$ws.Range("K116").Value = 6020.9287
$ws.Range("L116").Value = 14700
$ws.Range("M116").Value = -3726.9287
$ws.Range("N116").Value = -19288

$ws.Range("H122").Value = 3213.0205
$ws.Range("I122").Value = 2556
$ws.Range("K122").Value = 7668
$ws.Range("M122").Value = -5218

$ws.Range("H136").Value = 5205.3716
$ws.Range("I136").Value = 4765.3105
$ws.Range("K136").Value = 14295.9315
$ws.Range("M136").Value = -11745.9315

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8304.895
$ws.Range("I3").Value = 6020.9287
$ws.Range("J3").Value = 14700
$ws.Range("K3").Value = 6020.9287
$ws.Range("L3").Value = 14700
$ws.Range("M3").Value = -5906.9287
$ws.Range("N3").Value = -14928

$ws.Range("H64").Value = 3902.3333
$ws.Range("I64").Value = 850
$ws.Range("K64").Value = 850
$ws.Range("M64").Value = -625

$ws.Range("H67").Value = 3902.3333
$ws.Range("I67").Value = 850
$ws.Range("K67").Value = 850
$ws.Range("M67").Value = -70

$ws.Range("H105").Value = 13779.889
$ws.Range("I105").Value = 9299.6
$ws.Range("K105").Value = 9299.6
$ws.Range("M105").Value = -7552.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7424.353
$ws.Range("J31").Value = 9226.5
$ws.Range("L31").Value = 9226.5
$ws.Range("N31").Value = -9816.5

$ws.Range("H34").Value = 7424.353
$ws.Range("J34").Value = 9226.5
$ws.Range("L34").Value = 9226.5
$ws.Range("N34").Value = -9630.5

$ws.Range("H58").Value = 40008100
$ws.Range("I58").Value = 71434750
$ws.Range("J58").Value = 10540
$ws.Range("K58").Value = 71434750
$ws.Range("L58").Value = 10540
$ws.Range("M58").Value = -71434547
$ws.Range("N58").Value = -10946

$ws.Range("H94").Value = 52636144
$ws.Range("I94").Value = 166670200
$ws.Range("J94").Value = 5034.5386
$ws.Range("K94").Value = 166670200
$ws.Range("L94").Value = 5034.5386
$ws.Range("M94").Value = -166669749
$ws.Range("N94").Value = -5936.5386

$ws.Range("H105").Value = 58828160
$ws.Range("I105").Value = 142858600
$ws.Range("K105").Value = 142858600
$ws.Range("M105").Value = -142856853

$ws.Range("H122").Value = 1508.1724
$ws.Range("I122").Value = 1149.1364
$ws.Range("J122").Value = 2636.5715
$ws.Range("K122").Value = 3447.4092
$ws.Range("L122").Value = 7909.7145
$ws.Range("M122").Value = -997.4092000000001
$ws.Range("N122").Value = -12809.7145

$ws.Range("H132").Value = 4690.2856
$ws.Range("I132").Value = 4579.522
$ws.Range("K132").Value = 13738.566
$ws.Range("M132").Value = -11208.566

$ws.Range("H134").Value = 34490224
$ws.Range("I134").Value = 142866500
$ws.Range("J134").Value = 6861.091
$ws.Range("K134").Value = 428599500
$ws.Range("L134").Value = 20583.273
$ws.Range("M134").Value = -428596965
$ws.Range("N134").Value = -25653.273

$ws.Range("H136").Value = 40008100
$ws.Range("I136").Value = 71434750
$ws.Range("J136").Value = 10540
$ws.Range("K136").Value = 214304250
$ws.Range("L136").Value = 31620
$ws.Range("M136").Value = -214301700
$ws.Range("N136").Value = -36720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 135947.3
$ws.Range("I5").Value = 691.625
$ws.Range("J5").Value = 385650.06
$ws.Range("K5").Value = 2074.875
$ws.Range("L5").Value = 1156950.18
$ws.Range("M5").Value = -1962.875
$ws.Range("N5").Value = -1157174.18

$ws.Range("H39").Value = 6900
$ws.Range("J39").Value = 6900
$ws.Range("L39").Value = 20700
$ws.Range("N39").Value = -21288

$ws.Range("H55").Value = 15856.125
$ws.Range("J55").Value = 19666.666
$ws.Range("L55").Value = 58999.99800000001
$ws.Range("N55").Value = -59353.99800000001

$ws.Range("H132").Value = 2039
$ws.Range("I132").Value = 1239.8
$ws.Range("K132").Value = 11158.2
$ws.Range("M132").Value = -8628.199999999999

$ws.Range("H134").Value = 9276.362999999999
$ws.Range("I134").Value = 9804
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 29412
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -24342
$ws.Range("N134").Value = -22140

$ws.Range("H135").Value = 135947.3
$ws.Range("I135").Value = 691.625
$ws.Range("J135").Value = 385650.06
$ws.Range("K135").Value = 6224.625
$ws.Range("L135").Value = 3470850.54
$ws.Range("M135").Value = -3689.625
$ws.Range("N135").Value = -3475920.54

$ws.Range("H139").Value = 102888.78
$ws.Range("J139").Value = 451750
$ws.Range("L139").Value = 1355250
$ws.Range("N139").Value = -1365530

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 34485680
$ws.Range("I132").Value = 55557050
$ws.Range("K132").Value = 166671150
$ws.Range("M132").Value = -166668620

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5237
$ws.Range("I7").Value = 4464.5
$ws.Range("J7").Value = 5546
$ws.Range("K7").Value = 4464.5
$ws.Range("L7").Value = 5546
$ws.Range("M7").Value = -4352.5
$ws.Range("N7").Value = -5770

$ws.Range("H22").Value = 4519.6
$ws.Range("I22").Value = 2878.7896
$ws.Range("J22").Value = 9715.5
$ws.Range("K22").Value = 2878.7896
$ws.Range("L22").Value = 9715.5
$ws.Range("M22").Value = -2583.7896
$ws.Range("N22").Value = -10305.5

$ws.Range("H27").Value = 4519.6
$ws.Range("I27").Value = 2878.7896
$ws.Range("J27").Value = 9715.5
$ws.Range("K27").Value = 2878.7896
$ws.Range("L27").Value = 9715.5
$ws.Range("M27").Value = -2771.7896
$ws.Range("N27").Value = -9929.5

$ws.Range("H40").Value = 6956.143
$ws.Range("I40").Value = 4673.25
$ws.Range("K40").Value = 4673.25
$ws.Range("M40").Value = -4537.25

$ws.Range("H122").Value = 5152.773
$ws.Range("I122").Value = 4918.55
$ws.Range("K122").Value = 14755.65
$ws.Range("M122").Value = -12305.65

$ws.Range("H126").Value = 5237
$ws.Range("I126").Value = 4464.5
$ws.Range("J126").Value = 5546
$ws.Range("K126").Value = 13393.5
$ws.Range("L126").Value = 16638
$ws.Range("M126").Value = -10923.5
$ws.Range("N126").Value = -21578

$ws.Range("H136").Value = 107206110
$ws.Range("I136").Value = 50086708
$ws.Range("J136").Value = 250004600
$ws.Range("K136").Value = 150260124
$ws.Range("L136").Value = 750013800
$ws.Range("M136").Value = -150257574
$ws.Range("N136").Value = -750018900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 41969.2
$ws.Range("J74").Value = 37462.75
$ws.Range("L74").Value = 37462.75
$ws.Range("N74").Value = -39334.75

$ws.Range("H77").Value = 41969.2
$ws.Range("J77").Value = 37462.75
$ws.Range("L77").Value = 112388.25
$ws.Range("N77").Value = -121748.25

$ws.Range("H107").Value = 1225.5
$ws.Range("I107").Value = 1374.9131
$ws.Range("K107").Value = 4124.7393
$ws.Range("M107").Value = -2204.7393

$ws.Range("H126").Value = 5450.125
$ws.Range("J126").Value = 7119.4
$ws.Range("L126").Value = 21358.2
$ws.Range("N126").Value = -26298.2

$ws.Range("H132").Value = 8031.974
$ws.Range("I132").Value = 7705.125
$ws.Range("J132").Value = 8554.933999999999
$ws.Range("K132").Value = 23115.375
$ws.Range("L132").Value = 25664.802
$ws.Range("M132").Value = -20585.375
$ws.Range("N132").Value = -30724.802
